# Loan RBI, Variable Instalments
#
# On the "Repayment Schedule" sheet, insert a new (blank) column before
# column N - this splits the old "In Advance"/"Outstanding" block so a
# fresh variable-instalment column can be added later, shifting the old
# N/O/P columns to O/P/Q. Then make "Repayment Schedule" the active sheet
# (it was "Transactions" before) with R7 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before column N (shifts N->O, O-> (new, blank),
# P->Q, carrying each column's formatting along with it).
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with R7 selected -
# this also clears the previous tab-selected state on "Transactions".
$ws.Activate()
$ws.Range("R7").Select() | Out-Null
